$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.579.78"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "2.584.42"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.28"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.50"
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("E8").Value = "  -6.70%  "
$ws.Range("D9").Value = "2.588.09"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.60"
$ws.Range("E10").Value = "  +7.38%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "3.041.04"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("D15").Value = "60.523.42"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.50"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("D18").Value = "2.586.06"
$ws.Range("E18").Value = "  -2.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.79"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.17"
$ws.Range("E20").Value = "  +3.84%  "
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("E28").Value = "  +2.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.31"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.57"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.70"
$ws.Range("E34").Value = "  +3.72%  "
$ws.Range("E35").Value = "  +2.65%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.858"
$ws.Range("E37").Value = "  +10.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.851"
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.47"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.76"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.94"
$ws.Range("E41").Value = "  +2.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "296.21"
$ws.Range("E42").Value = "  +2.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0996"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0557"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.996"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.76"
$ws.Range("E47").Value = "  +3.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.86"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").Value = "2.001.16"
$ws.Range("E51").Value = "  +0.03%  "
